$d = $word.ActiveDocument

# --- Collapse the source paragraphs -----------------------------------------
# The original document told the story across three text paragraphs, a blank
# paragraph, and a trailing "Was fuer ein Index?" question paragraph. The new
# text is a single rewritten paragraph, so first merge everything from the
# end of paragraph 1 through the end of the last paragraph away, leaving one
# empty paragraph (paragraph 1) in place (its own paragraph mark / rsids are
# preserved since we never touch that paragraph's own range boundary).
$p1 = $d.Paragraphs.Item(1)
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
if ($lastPara.Range.End -gt $p1.Range.End) {
    $tailRange = $d.Range($p1.Range.End, $lastPara.Range.End)
    $tailRange.Delete()
}

# --- New text, chunked the same way the author's runs were ------------------
$chunks = @(
    'Zuerst wird dem Programm eine Datei übergeben',
    ', die für das schreiben und lesen ',
    'geöffnet wird ',
    '.',
    ' Um ',
    'anschließend ',
    'die übergebene Datei schrittweise Auslesen zu können wird ein Index des gelesenen Arrays, ein Index des geschriebenen Arrays und ein Zeiger für die Arrays angelegt. ',
    'Im nächsten Schritt wird überprüft, ob es sich beim Inhalt der Datei um Zahlen handelt. Besitzt die Datei keine Zahl, so wird sie darauf gleich geschlossen. Sind in der Datei Zahlen vorhanden beginnt der Hauptprozess. Die Zahl wird am Index des gelesenes Arrays gelesen. Anschließend wird die ',
    'Zahl in einer Variable zwischen',
    'gespeichert. Der Wert der Variable wird daraufhin überprüft, ob es sich zum einen um die Zahl -1 handelt, zum anderen nicht um die Zahl -1. ',
    'Wird nun festgestellt, dass es sich bei dem Wert der Variable um die Zahl -1 handelt, so wird anschließend die -1 im Index des geschriebenen Arrays geschrieben, der Index wird ',
    'um ein',
    's erhöht und daraufhin wird der Wert der Variable im Index geschrieben und zum Abschluss',
    ' wieder der Index um ein',
    's erhöht.',
    ' Der Index ',
    'des gelesenen Arrays wird daraufhin erhöht, sowohl auch der Index des geschriebenen Arrays. ',
    'E',
    'ine weitere Zahl wird vom Index des gelesenen Arrays gelesen. ',
    'Handelt es sich nun bei einem weiteren Durchlauf um eine andere Zahl außer -1, wird der Zeiger',
    ' um ein',
    's erhöht und die Zahl an der entsprechenden Position ausgelesen. Nun werden wieder zwei Fälle betrachtet, zum einen ob die gespeichert Zahl gleich der Zahl an der Zeigerposition ist, zum anderen ob die Zahl ungleich der Zahl an der Zeigerposition ist.',
    ' Bei ungleicher Zahl wird der Wert direkt in eine neue Datei an der Position des Index des geschriebenes Arrays gelesen',
    ', der Index des gelesenes Arrays und der Index des geschriebenes Arrays ',
    'werden jeweils um ein',
    's erhöht.',
    ' ',
    'E',
    'ine weitere Zahl wird vom Index des gelesenen Arrays gelesen.',
    ' Beim letzten Fall ist der gespeicherte Wert und die Zahl an der Position des Zeigers gleich',
    '. Zuerst wird im Index des geschriebenen Arrays die -1 geschrieben, anschließend wird der Index erhöht. Daraufhin wird die Zahl im Index des geschriebenen Arrays geschrieben u',
    'nd der Index wird wieder um ein',
    's erhöht. Danach wird die Anzahl der Bytes ermittelt und diese dann am Index des ge',
    'schriebenes Arrays geschrieben, der Index wird wieder um eines erhöht und der Index des gelesenen Arrays wird die Anzahl um eins erhöht.',
    ' Ist die Datei nicht vollständig ausgelesen, wird wieder eine weitere Zahl am Index des gelesenen Arrays gelesen und der Prozess beginnt von Neuem, bis die Datei endgültig gelesen worden ist.'
)

# --- Write the first chunk over the whole (now single) paragraph, then append
#     every other chunk right after the paragraph's current end. ---
$p1 = $d.Paragraphs.Item(1)
$firstRange = $d.Range($p1.Range.Start, $p1.Range.End)
$firstRange.Text = $chunks[0]

for ($i = 1; $i -lt $chunks.Length; $i++) {
    $endPos = $d.Paragraphs.Item(1).Range.End
    $insertionPoint = $d.Range($endPos, $endPos)
    $insertionPoint.InsertAfter($chunks[$i])
}
